# "Finally got rid of health bar bug!"
# Remove the "Career Criminal" / "Criminal-damage-arson" row (row 7) from
# the crime-stats table entirely (it was a duplicate/bugged entry), which
# shifts the rows below it up by one (Professional -> row 7, Gankster ->
# row 8, and the city list below moves from rows 14-23 up to rows 13-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 7 ("Career Criminal" / "Criminal-damage-arson"),
# shifting everything below it up by one row.
$ws.Rows(7).Delete() | Out-Null

# Match the new selection left behind in the saved file.
$ws.Range("B8").Select() | Out-Null
